# Refresh scraped crypto Price (D) and Volume(1h) (E) columns with
# the latest run's values. The source writes these as plain text
# (matches the pre-existing inline-string cell type), so the column
# is pre-set to Text format before assigning the numeric-looking
# Price entries -- otherwise Excel auto-converts "239.74" etc. to a
# number. The format is reverted to Normal afterwards so no stray
# number formatting is left on the cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = "43.904.24"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.348.19"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "239.74"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("D7").Value = "73.52"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "61.33"
$ws.Range("E11").Value = "  +7.29%  "
$ws.Range("D12").Value = "33.48"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "7.27"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "16.16"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "2.343.76"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "43.761.63"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").Value = "77.82"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").Value = "252.57"
$ws.Range("D23").Value = "3.81"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "1.84"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "10.41"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "175.87"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "22.19"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -5.11%  "
$ws.Range("D36").Value = "3.77"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").Value = "6.42"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "5.38"
$ws.Range("E40").Value = "  +12.06%  "
$ws.Range("D41").Value = "65.77"
$ws.Range("E41").Value = "  +15.72%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "9.14"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").Value = "0.200"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "97.99"
$ws.Range("E50").Value = "  -3.09%  "

$ws.Range("D2:D50").Style = "Normal"
